# Update the simulated-game transition-probability matrix on Sheet1.
# More games were simulated, so the per-row game counts (and therefore the
# row-normalized probabilities in columns B:S) changed for several states.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "B2" = 0.1838006230529595
    "C2" = 0.5825545171339563
    "J2" = 0.01557632398753894
    "O2" = 0.003115264797507788
    "P2" = 0.1308411214953271
    "S2" = 0.08411214953271028

    "B3" = 0.01041666666666667
    "C3" = 0.02083333333333333
    "J3" = 0.03645833333333334
    "P3" = 0.7604166666666666
    "S3" = 0.171875

    "J4" = 0.02222222222222222
    "P4" = 0.6666666666666666
    "S4" = 0.3111111111111111

    "B6" = 0.06349206349206349
    "D6" = 0.01587301587301587
    "F6" = 0.1216931216931217
    "J6" = 0.1534391534391534
    "O6" = 0.04761904761904762
    "Q6" = 0.164021164021164
    "R6" = 0.06878306878306878
    "S6" = 0.3650793650793651

    "B7" = 0.1222707423580786
    "D7" = 0.01746724890829694
    "F7" = 0.03056768558951965
    "J7" = 0.1310043668122271
    "O7" = 0.008733624454148471
    "Q7" = 0.1266375545851528
    "R7" = 0.04803493449781659
    "S7" = 0.5152838427947598

    "B8" = 0.08968609865470852
    "D8" = 0.02017937219730942
    "F8" = 0.06502242152466367
    "J8" = 0.1233183856502242
    "O8" = 0.02017937219730942
    "Q8" = 0.1345291479820628
    "R8" = 0.08295964125560538
    "S8" = 0.4641255605381166

    "B9" = 0.1092896174863388
    "D9" = 0.01092896174863388
    "F9" = 0.09836065573770492
    "J9" = 0.1693989071038251
    "O9" = 0.00546448087431694
    "Q9" = 0.1366120218579235
    "R9" = 0.08743169398907104
    "S9" = 0.3825136612021858

    "B10" = 0.1381461675579323
    "D10" = 0.02584670231729055
    "F10" = 0.05793226381461675
    "J10" = 0.1631016042780749
    "O10" = 0.0196078431372549
    "Q10" = 0.1755793226381462
    "R10" = 0.0659536541889483
    "S10" = 0.3538324420677362

    "G11" = 0.1354838709677419
    "J11" = 0.06129032258064516
    "K11" = 0.1774193548387097
    "L11" = 0.6064516129032258
    "S11" = 0.01935483870967742

    "G12" = 0.7889447236180904
    "J12" = 0.1457286432160804
    "K12" = 0.005025125628140704
    "L12" = 0.03015075376884422
    "S12" = 0.03015075376884422

    "G13" = 0.6724137931034483
    "J13" = 0.2758620689655172
    "S13" = 0.05172413793103448

    "J14" = 0.25
    "S14" = 0.25

    "F15" = 0.0198019801980198
    "H15" = 0.1831683168316832
    "I15" = 0.08415841584158416
    "J15" = 0.2475247524752475
    "K15" = 0.09900990099009901
    "M15" = 0.03465346534653466
    "O15" = 0.08415841584158416
    "S15" = 0.2475247524752475

    "F16" = 0.01923076923076923
    "H16" = 0.2067307692307692
    "I16" = 0.04807692307692308
    "J16" = 0.3317307692307692
    "K16" = 0.1442307692307692
    "M16" = 0.03365384615384615
    "O16" = 0.07692307692307693
    "S16" = 0.1394230769230769

    "F17" = 0.01166180758017493
    "H17" = 0.1749271137026239
    "I17" = 0.09329446064139942
    "J17" = 0.4518950437317784
    "K17" = 0.08746355685131195
    "M17" = 0.02623906705539359
    "N17" = 0.002915451895043732
    "O17" = 0.04664723032069971
    "S17" = 0.1049562682215743

    "F18" = 0.01948051948051948
    "H18" = 0.1753246753246753
    "I18" = 0.08441558441558442
    "J18" = 0.4025974025974026
    "K18" = 0.09090909090909091
    "M18" = 0.006493506493506494
    "N18" = 0.006493506493506494
    "O18" = 0.07142857142857142
    "S18" = 0.1428571428571428

    "F19" = 0.01215559157212318
    "H19" = 0.2204213938411669
    "I19" = 0.08670988654781199
    "J19" = 0.3290113452188007
    "K19" = 0.1312803889789303
    "M19" = 0.02836304700162074
    "N19" = 0.002431118314424636
    "O19" = 0.05672609400324149
    "S19" = 0.1329011345218801
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
